# Trade #61 closed at 2026-02-16 21:34:08 - leadlag DOWN +0.000%
#
# This script applies the following changes to live_trading_results.xlsx:
#  1. A new "leadlag" trade (#61) was opened (row 50 on the "leadlag" sheet).
#  2. The previously-open "momentum" trade #39 (row 11 on the "momentum"
#     sheet) closed, so its exit fields are filled in and a mirrored row is
#     appended to the "All Trades" sheet (row 40).
#  3. The Summary and Comparison sheets are refreshed with the updated
#     aggregate statistics that result from the above.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing "date-looking" strings
# (e.g. 2026-02-16) to stay plain text instead of being auto-converted
# into an Excel date serial number.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
}

# ---------------------------------------------------------------------
# 1. "leadlag" sheet - append trade #61 (newly opened) as row 50
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Range("A50").Value = 61
Set-TextValue $leadlag.Range("B50") "2026-02-16"
$leadlag.Range("C50").Value = "21:34:08"
$leadlag.Range("D50").Value = "leadlag"
$leadlag.Range("E50").Value = "DOWN"
$leadlag.Range("F50").Value = 68763.33
$leadlag.Range("G50").Value = ""
$leadlag.Range("H50").Value = "OPEN"
$leadlag.Range("I50").Value = 0
$leadlag.Range("J50").Value = 0
$leadlag.Range("K50").Value = 0.75
$leadlag.Range("L50").Value = "Coinbase leading with -0.090% move"
$leadlag.Range("M50").Value = ""
$leadlag.Range("N50").Value = 0

# ---------------------------------------------------------------------
# 2. "momentum" sheet - close out trade #39 (row 11)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Range("G11").Value = 68563.063352
$momentum.Range("H11").Value = "CLOSED"
$momentum.Range("I11").Value = 0.1184
$momentum.Range("J11").Value = 1.18
$momentum.Range("M11").Value = "time_exit_5min"
$momentum.Range("N11").Value = 5

# ---------------------------------------------------------------------
# 3. "All Trades" sheet - append the mirrored, now-closed row for trade #39
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A40").Value = 39
Set-TextValue $allTrades.Range("B40") "2026-02-16"
$allTrades.Range("C40").Value = "21:29:05"
$allTrades.Range("D40").Value = "momentum"
$allTrades.Range("E40").Value = "DOWN"
$allTrades.Range("F40").Value = 68644.355
$allTrades.Range("G40").Value = 68563.063352
$allTrades.Range("H40").Value = "CLOSED"
$allTrades.Range("I40").Value = 0.1184
$allTrades.Range("J40").Value = 1.18
$allTrades.Range("K40").Value = 0.9
$allTrades.Range("L40").Value = "Downward momentum: -0.417% over 10 samples"
$allTrades.Range("M40").Value = "time_exit_5min"
$allTrades.Range("N40").Value = 5

# ---------------------------------------------------------------------
# 4. "Summary" sheet - refreshed aggregate stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 39
Set-TextValue $summary.Range("D2") "64.1%"
Set-TextValue $summary.Range("E2") "+7.6677%"
Set-TextValue $summary.Range("F2") "+0.1966%"

Set-TextValue $summary.Range("D4") "66.7%"
Set-TextValue $summary.Range("E4") "+3.0385%"
Set-TextValue $summary.Range("F4") "+0.2532%"

# ---------------------------------------------------------------------
# 5. "Comparison" sheet - refreshed momentum row stats
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

Set-TextValue $comparison.Range("C3") "66.7%"
Set-TextValue $comparison.Range("D3") "3.70"
Set-TextValue $comparison.Range("E3") "+0.5204%"
Set-TextValue $comparison.Range("G3") "0.93"
